$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value = 0.027123
$ws.Range("H2").Value = 0.081369
$ws.Range("I2").Value = 0.0960827240265261
$ws.Range("J2").Value = 0.09608272402652611
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.5
$ws.Range("M2").Value = 0.016657
$ws.Range("N2").Value = 0.033314
$ws.Range("Q2").Value = 0.000451787811
$ws.Range("R2").Value = 0.002710726866
$ws.Range("S2").Value = 0.0960827240265261
$ws.Range("T2").Value = 0.09608272402652611

# Row 3 updates
$ws.Range("I3").Value = 0.9039172759734738
$ws.Range("J3").Value = 0.9039172759734738
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.5
$ws.Range("M3").Value = 0.016657
$ws.Range("N3").Value = 0.033314
$ws.Range("Q3").Value = 0.004250283405
$ws.Range("R3").Value = 0.02550170043
$ws.Range("S3").Value = 0.9039172759734738
$ws.Range("T3").Value = 0.9039172759734738
